$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10492.363
$ws.Range("I32").Value = 1829
$ws.Range("J32").Value = 13741.125
$ws.Range("K32").Value = 1829
$ws.Range("L32").Value = 13741.125
$ws.Range("M32").Value = -1503
$ws.Range("N32").Value = -14393.125

$ws.Range("H33").Value = 527.9286
$ws.Range("I33").Value = 539.3
$ws.Range("K33").Value = 539.3
$ws.Range("M33").Value = -310.3

$ws.Range("H40").Value = 5555.4443
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 6833.1665
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 6833.1665
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -7183.1665

$ws.Range("H43").Value = 8125.8945
$ws.Range("I43").Value = 6125
$ws.Range("K43").Value = 6125
$ws.Range("M43").Value = -6056

$ws.Range("H107").Value = 2901.6316
$ws.Range("I107").Value = 1514.5625
$ws.Range("J107").Value = 10299.333
$ws.Range("K107").Value = 1514.5625
$ws.Range("L107").Value = 10299.333
$ws.Range("M107").Value = 405.4375
$ws.Range("N107").Value = -14139.333

$ws.Range("H137").Value = 7614.6333
$ws.Range("J137").Value = 9931.691999999999
$ws.Range("L137").Value = 29795.076
$ws.Range("N137").Value = -34895.076

$ws.Range("H138").Value = 3105.137
$ws.Range("I138").Value = 981.4211
$ws.Range("J138").Value = 5410.8857
$ws.Range("K138").Value = 2944.2633
$ws.Range("L138").Value = 16232.6571
$ws.Range("M138").Value = 2195.7367
$ws.Range("N138").Value = -26512.6571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1794444.1
$ws.Range("I32").Value = 819.9
$ws.Range("J32").Value = 9267878
$ws.Range("K32").Value = 819.9
$ws.Range("L32").Value = 9267878
$ws.Range("M32").Value = -532.9
$ws.Range("N32").Value = -9268452

$ws.Range("H74").Value = 2474.2812
$ws.Range("I74").Value = 2805.92
$ws.Range("J74").Value = 1289.8572
$ws.Range("K74").Value = 2805.92
$ws.Range("L74").Value = 1289.8572
$ws.Range("M74").Value = -1931.92
$ws.Range("N74").Value = -3037.8572

$ws.Range("H77").Value = 2474.2812
$ws.Range("I77").Value = 2805.92
$ws.Range("J77").Value = 1289.8572
$ws.Range("K77").Value = 14029.6
$ws.Range("L77").Value = 6449.286
$ws.Range("M77").Value = -9661.6
$ws.Range("N77").Value = -15185.286

$ws.Range("H102").Value = 5288.722
$ws.Range("I102").Value = 5512.5625
$ws.Range("K102").Value = 5512.5625
$ws.Range("M102").Value = -3890.5625

$ws.Range("H122").Value = 2425.1396
$ws.Range("I122").Value = 1664.2
$ws.Range("K122").Value = 4992.6
$ws.Range("M122").Value = -2542.6

$ws.Range("H132").Value = 896651.9399999999
$ws.Range("I132").Value = 1011768.94
$ws.Range("J132").Value = 205949.8
$ws.Range("K132").Value = 3035306.82
$ws.Range("L132").Value = 617849.3999999999
$ws.Range("M132").Value = -3032776.82
$ws.Range("N132").Value = -622909.3999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 3162.6
$ws.Range("J36").Value = 6990
$ws.Range("L36").Value = 6990
$ws.Range("N36").Value = -8058

$ws.Range("H82").Value = 14460.667
$ws.Range("I82").Value = 14460.667
$ws.Range("K82").Value = 14460.667
$ws.Range("M82").Value = -14077.667

$ws.Range("H85").Value = 14460.667
$ws.Range("I85").Value = 14460.667
$ws.Range("K85").Value = 14460.667
$ws.Range("M85").Value = -13134.667

$ws.Range("H99").Value = 7290.3623
$ws.Range("J99").Value = 8847.174000000001
$ws.Range("L99").Value = 8847.174000000001
$ws.Range("N99").Value = -11843.174

$ws.Range("H107").Value = 7694614
$ws.Range("I107").Value = 8335261.5
$ws.Range("K107").Value = 8335261.5
$ws.Range("M107").Value = -8333341.5

$ws.Range("H134").Value = 870246.4
$ws.Range("I134").Value = 995777.0600000001
$ws.Range("K134").Value = 2987331.18
$ws.Range("M134").Value = -2984796.18

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4527.8184
$ws.Range("I31").Value = 1622
$ws.Range("K31").Value = 1622
$ws.Range("M31").Value = -1327

$ws.Range("H34").Value = 4527.8184
$ws.Range("I34").Value = 1622
$ws.Range("K34").Value = 1622
$ws.Range("M34").Value = -1420

$ws.Range("H38").Value = 15406.2
$ws.Range("I38").Value = 18983
$ws.Range("J38").Value = 10041
$ws.Range("K38").Value = 18983
$ws.Range("L38").Value = 10041
$ws.Range("M38").Value = -18606
$ws.Range("N38").Value = -10795

$ws.Range("H46").Value = 15406.2
$ws.Range("I46").Value = 18983
$ws.Range("J46").Value = 10041
$ws.Range("K46").Value = 18983
$ws.Range("L46").Value = 10041
$ws.Range("M46").Value = -18772
$ws.Range("N46").Value = -10463

$ws.Range("H107").Value = 703.7826
$ws.Range("J107").Value = 687.4167
$ws.Range("L107").Value = 687.4167
$ws.Range("N107").Value = -4527.4167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1439.6154
$ws.Range("I18").Value = 1023.7143
$ws.Range("K18").Value = 3071.1429
$ws.Range("M18").Value = -2902.1429

$ws.Range("H29").Value = 999
$ws.Range("I29").Value = 999
$ws.Range("K29").Value = 2997
$ws.Range("M29").Value = -2720

$ws.Range("H39").Value = 4519.6665
$ws.Range("J39").Value = 4519.6665
$ws.Range("L39").Value = 13558.9995
$ws.Range("N39").Value = -14146.9995

$ws.Range("H55").Value = 3934.4
$ws.Range("I55").Value = 3825
$ws.Range("J55").Value = 3974.182
$ws.Range("K55").Value = 11475
$ws.Range("L55").Value = 11922.546
$ws.Range("M55").Value = -11298
$ws.Range("N55").Value = -12276.546

$ws.Range("H64").Value = 13789.1
$ws.Range("I64").Value = 9998
$ws.Range("J64").Value = 14736.875
$ws.Range("K64").Value = 29994
$ws.Range("L64").Value = 44210.625
$ws.Range("M64").Value = -29724
$ws.Range("N64").Value = -44750.625

$ws.Range("H67").Value = 13789.1
$ws.Range("I67").Value = 9998
$ws.Range("J67").Value = 14736.875
$ws.Range("K67").Value = 29994
$ws.Range("L67").Value = 44210.625
$ws.Range("M67").Value = -29058
$ws.Range("N67").Value = -46082.625

$ws.Range("H87").Value = 16416.5
$ws.Range("I87").Value = 11963.363
$ws.Range("J87").Value = 23414.285
$ws.Range("K87").Value = 35890.089
$ws.Range("L87").Value = 70242.855
$ws.Range("M87").Value = -34642.089
$ws.Range("N87").Value = -72738.855

$ws.Range("H90").Value = 16416.5
$ws.Range("I90").Value = 11963.363
$ws.Range("J90").Value = 23414.285
$ws.Range("K90").Value = 107670.267
$ws.Range("L90").Value = 210728.565
$ws.Range("M90").Value = -101430.267
$ws.Range("N90").Value = -223208.565

$ws.Range("H94").Value = 7845.174
$ws.Range("J94").Value = 8050.6313
$ws.Range("L94").Value = 24151.8939
$ws.Range("N94").Value = -25503.8939

$ws.Range("H132").Value = 3403.7
$ws.Range("I132").Value = 1510.8
$ws.Range("K132").Value = 13597.2
$ws.Range("M132").Value = -11067.2

$ws.Range("H140").Value = 47103660
$ws.Range("I140").Value = 63726644
$ws.Range("J140").Value = 5198.8335
$ws.Range("K140").Value = 191179932
$ws.Range("L140").Value = 15596.5005
$ws.Range("M140").Value = -191174752
$ws.Range("N140").Value = -25956.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 9500
$ws.Range("J49").Value = 9500
$ws.Range("L49").Value = 9500
$ws.Range("N49").Value = -9868

$ws.Range("H102").Value = 4568.7456
$ws.Range("I102").Value = 3592.6562
$ws.Range("J102").Value = 5725.593
$ws.Range("K102").Value = 3592.6562
$ws.Range("L102").Value = 5725.593
$ws.Range("M102").Value = -1970.6562
$ws.Range("N102").Value = -8969.593000000001

$ws.Range("H122").Value = 5806.2085
$ws.Range("I122").Value = 5108.4614
$ws.Range("J122").Value = 6630.8184
$ws.Range("K122").Value = 15325.3842
$ws.Range("L122").Value = 19892.4552
$ws.Range("M122").Value = -12875.3842
$ws.Range("N122").Value = -24792.4552

$ws.Range("H126").Value = 27783588
$ws.Range("I126").Value = 62501904
$ws.Range("J126").Value = 8935.6
$ws.Range("K126").Value = 187505712
$ws.Range("L126").Value = 26806.8
$ws.Range("M126").Value = -187503242
$ws.Range("N126").Value = -31746.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 3333.3333
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5224

$ws.Range("H12").Value = 14286857
$ws.Range("J12").Value = 1624.5
$ws.Range("L12").Value = 1624.5
$ws.Range("N12").Value = -1964.5

$ws.Range("H15").Value = 3333.3333
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5340

$ws.Range("H46").Value = 50002596
$ws.Range("J46").Value = 71431770
$ws.Range("L46").Value = 71431770
$ws.Range("N46").Value = -71432146

$ws.Range("H55").Value = 3525.027
$ws.Range("I55").Value = 2204.1853
$ws.Range("J55").Value = 7091.3
$ws.Range("K55").Value = 2204.1853
$ws.Range("L55").Value = 7091.3
$ws.Range("M55").Value = -2031.1853
$ws.Range("N55").Value = -7437.3

$ws.Range("H136").Value = 107150920
$ws.Range("I136").Value = 50008292
$ws.Range("K136").Value = 150024876
$ws.Range("M136").Value = -150022326

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 12078.857
$ws.Range("I2").Value = 12078.857
$ws.Range("K2").Value = 12078.857
$ws.Range("M2").Value = -11966.857

$ws.Range("H107").Value = 6667243.5
$ws.Range("I107").Value = 9524375
$ws.Range("K107").Value = 28573125
$ws.Range("M107").Value = -28571205

$ws.Range("H126").Value = 4859.2383
$ws.Range("J126").Value = 7908.8
$ws.Range("L126").Value = 23726.4
$ws.Range("N126").Value = -28666.4
